# Update the date header and the 25 division problems in the practice-sheet
# table. Order matters: row 12/col 1 turns "97÷6=" into "66÷5=", and
# "66÷5=" is also the *original* text of row 0/col 0 (which is replaced
# first, earlier in document order), so processing top-to-bottom avoids any
# accidental re-matching.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-01 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-02 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("66÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷9=", 2) | Out-Null
$d.Content.Find.Execute("78÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷2=", 2) | Out-Null
$d.Content.Find.Execute("24÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷6=", 2) | Out-Null
$d.Content.Find.Execute("40÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=", 2) | Out-Null
$d.Content.Find.Execute("53÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷8=", 2) | Out-Null
$d.Content.Find.Execute("99÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷9=", 2) | Out-Null
$d.Content.Find.Execute("61÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷7=", 2) | Out-Null
$d.Content.Find.Execute("87÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷4=", 2) | Out-Null
$d.Content.Find.Execute("16÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=", 2) | Out-Null
$d.Content.Find.Execute("26÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=", 2) | Out-Null
$d.Content.Find.Execute("79÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷5=", 2) | Out-Null
$d.Content.Find.Execute("16÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷4=", 2) | Out-Null
$d.Content.Find.Execute("77÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷8=", 2) | Out-Null
$d.Content.Find.Execute("65÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷9=", 2) | Out-Null
$d.Content.Find.Execute("28÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷3=", 2) | Out-Null
$d.Content.Find.Execute("60÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷3=", 2) | Out-Null
$d.Content.Find.Execute("97÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷5=", 2) | Out-Null
$d.Content.Find.Execute("37÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷2=", 2) | Out-Null
$d.Content.Find.Execute("68÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=", 2) | Out-Null
$d.Content.Find.Execute("30÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=", 2) | Out-Null
$d.Content.Find.Execute("51÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=", 2) | Out-Null
$d.Content.Find.Execute("77÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=", 2) | Out-Null
$d.Content.Find.Execute("47÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=", 2) | Out-Null
$d.Content.Find.Execute("18÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=", 2) | Out-Null
$d.Content.Find.Execute("73÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷2=", 2) | Out-Null
